# Fixed logic of incorrect SAN alert not followed by the SAN Input box
#
# Summary of the change:
#  - The "4.2 Items" sheet's running Desktop Mini G9 counters are bumped
#    up by the 4 SAN numbers that were actually entered correctly this
#    time (109 -> 114, 110 -> 115).
#  - The "4.2 Timestamps" log sheet gets the 4 missing "add" entries for
#    those SAN numbers appended (rows 34-37).
#  - The "All SANs" sheet gets the matching 4 rows filled in (rows 12-15,
#    which already existed as blank placeholder rows).

$wb = $excel.ActiveWorkbook

$itemsSheet = $wb.Worksheets.Item("4.2 Items")
$timestampsSheet = $wb.Worksheets.Item("4.2 Timestamps")
$allSansSheet = $wb.Worksheets.Item("All SANs")

# 1. Update the Desktop Mini G9 counters on "4.2 Items" (row 2).
$itemsSheet.Cells.Item(2, 2).Value = 114
$itemsSheet.Cells.Item(2, 3).Value = 115

# 2. Append the 4 new "add" log rows on "4.2 Timestamps" (columns:
#    Item, Action, SAN Number, Time).
$newEntries = @(
    @("2023-12-31 16:04:39", "Desktop Mini G9", "add", "SAN655443"),
    @("2023-12-31 16:05:19", "Desktop Mini G9", "add", "SAN434343"),
    @("2023-12-31 16:08:18", "Desktop Mini G9", "add", "SAN111111"),
    @("2023-12-31 16:08:38", "Desktop Mini G9", "add", "SAN111444")
)

$startRow = 34
for ($i = 0; $i -lt $newEntries.Count; $i++) {
    $row = $startRow + $i
    $entry = $newEntries[$i]
    $timestampsSheet.Cells.Item($row, 1).Value = $entry[0]
    $timestampsSheet.Cells.Item($row, 2).Value = $entry[1]
    $timestampsSheet.Cells.Item($row, 3).Value = $entry[2]
    $timestampsSheet.Cells.Item($row, 4).Value = $entry[3]
}

# 3. Fill in the matching rows on "All SANs" (columns: Item, SAN Number,
#    Time). Rows 12-15 already exist there as blank placeholder rows.
$sanEntries = @(
    @("Desktop Mini G9", "SAN655443", "2023-12-31 16:04:39"),
    @("Desktop Mini G9", "SAN434343", "2023-12-31 16:05:19"),
    @("Desktop Mini G9", "SAN111111", "2023-12-31 16:08:18"),
    @("Desktop Mini G9", "SAN111444", "2023-12-31 16:08:38")
)

$sanStartRow = 12
for ($i = 0; $i -lt $sanEntries.Count; $i++) {
    $row = $sanStartRow + $i
    $entry = $sanEntries[$i]
    $allSansSheet.Cells.Item($row, 1).Value = $entry[0]
    $allSansSheet.Cells.Item($row, 2).Value = $entry[1]
    $allSansSheet.Cells.Item($row, 3).Value = $entry[2]
}
